# "Generate Report for Handback" - refresh the handback-status report with
# newly regenerated timestamps (and an updated priority code for the
# zh-cn / de-de handoff row) exactly as a fresh CI run would produce.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G), rows 2 and 4
# share the same text, regenerate it to the later run time.
$wsOverview.Range("G2").Value = "2016-08-25 12:16:21"
$wsOverview.Range("G4").Value = "2016-08-25 12:16:21"

# zh-cn sheet: Priority flipped from human-translate (ht) to machine-translate (mt)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime (H) / Correspond Handback DateTime (K)
$wsZhCn.Range("H2").Value = "2016-08-25 12:16:16"
$wsZhCn.Range("H4").Value = "2016-08-25 12:16:16"
$wsZhCn.Range("K2").Value = "2016-08-25 12:16:32"
$wsZhCn.Range("K4").Value = "2016-08-25 12:16:32"

# de-de sheet: Priority flipped from human-translate (ht) to machine-translate (mt)
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# de-de sheet: Correspond Handoff Datetime (H) shares text with Overview!G
$wsDeDe.Range("H2").Value = "2016-08-25 12:16:21"
$wsDeDe.Range("H4").Value = "2016-08-25 12:16:21"

# de-de sheet: Correspond Handback DateTime (K)
$wsDeDe.Range("K2").Value = "2016-08-25 12:16:40"
$wsDeDe.Range("K4").Value = "2016-08-25 12:16:40"
